# Apply updated crypto price/volume data to sheet1 (cols D and E, rows 2-51)
# D-column values are forced to text ("@" number format) before assignment
# so that numeric-looking strings (e.g. "7.00", "0.0000261") are stored verbatim
# instead of being reinterpreted/re-serialized as floating point numbers
# (which would lose trailing zeros, use scientific notation, or introduce
# binary floating point rounding noise). The cell style is then restored to
# "Normal" so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.547.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.761.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.759.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.390.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.765.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.483.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("E23").Value = "  -3.05%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.908.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.717.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.303"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "388.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.21%  "
